$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; this shifts the existing
# PercActivations / PercSegmentAreas / RelativeCAMImportance /
# PercActivationsRescaled columns (B:E) one place to the right (C:F),
# and leaves the segment-name column (A) untouched.
$ws.Columns("B").Insert()

# New header for the inserted column; give it the same (bold, bordered,
# centered) header formatting used by the other header cells.
$ws.Range("B1").Value = "segments"
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the existing segment-name labels (currently still in column A,
# rows 2-20) into the newly inserted column B, then replace column A
# with a plain numeric index (0-based row counter).
$lastRow = 20
for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# The column insert copied the bold/bordered header look into B2:B20;
# reset those cells back to the plain/default style used by the data
# columns (C:F), matching the target layout.
$ws.Range("C2").Copy()
$ws.Range("B2:B20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Output "done"
